# Re-sort the analytics rows (2-24) chronologically by the "Version" timestamp column (L).
# The original data was grouped by Repository; after this edit it is ordered by date so the
# combined/"total" series can be charted as a single continuous time line (plot total graph).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2: User @ 2021-09-14-03:28
$ws.Range("A2").Value = 0.625
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.268125
$ws.Range("H2").Value = 0.5
$ws.Range("I2").Value = 0.7681249999999999
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = "User"
$ws.Range("L2").Value = "2021-09-14-03:28"
$ws.Range("M2").Value = "2021-09-14-03:28"

# row 3: Content @ 2021-09-14-03:28
$ws.Range("A3").Value = 0.7142857142857143
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.8928571428571429
$ws.Range("D3").Value = 0.8928571428571429
$ws.Range("E3").Value = 0.8928571428571429
$ws.Range("F3").Value = 0.8928571428571429
$ws.Range("G3").Value = 0.2651785714285714
$ws.Range("H3").Value = 0.4464285714285714
$ws.Range("I3").Value = 0.7116071428571429
$ws.Range("J3").Value = 629
$ws.Range("K3").Value = "Content"
$ws.Range("L3").Value = "2021-09-14-03:28"
$ws.Range("M3").Value = "2021-09-14-03:28"

# row 4: Mobile_App @ 2021-09-14-03:28
$ws.Range("A4").Value = 0.5555555555555556
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2566666666666667
$ws.Range("H4").Value = 0.5
$ws.Range("I4").Value = 0.7566666666666667
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = "Mobile_App"
$ws.Range("L4").Value = "2021-09-14-03:28"
$ws.Range("M4").Value = "2021-09-14-03:28"

# row 5: Files @ 2021-09-14-03:28
$ws.Range("A5").Value = 0.5714285714285714
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2592857142857143
$ws.Range("H5").Value = 0.5
$ws.Range("I5").Value = 0.7592857142857143
$ws.Range("J5").Value = 54
$ws.Range("K5").Value = "Files"
$ws.Range("L5").Value = "2021-09-14-03:28"
$ws.Range("M5").Value = "2021-09-14-03:28"

# row 6: User @ 2021-09-27-02.31
$ws.Range("A6").Value = 0.625
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.268125
$ws.Range("H6").Value = 0.5
$ws.Range("I6").Value = 0.7681249999999999
$ws.Range("J6").Value = 70
$ws.Range("K6").Value = "User"
$ws.Range("L6").Value = "2021-09-27-02.31"
$ws.Range("M6").Value = "2021-09-27-02.31"

# row 7: Mobile_App @ 2021-09-27-03.13
$ws.Range("A7").Value = 0.5555555555555556
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2566666666666667
$ws.Range("H7").Value = 0.5
$ws.Range("I7").Value = 0.7566666666666667
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = "Mobile_App"
$ws.Range("L7").Value = "2021-09-27-03.13"
$ws.Range("M7").Value = "2021-09-27-03.13"

# row 8: Content @ 2021-09-27-03.22
$ws.Range("A8").Value = 0.7115384615384616
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0.9423076923076923
$ws.Range("D8").Value = 0.9423076923076923
$ws.Range("E8").Value = 0.9423076923076923
$ws.Range("F8").Value = 0.9423076923076923
$ws.Range("G8").Value = 0.2728846153846154
$ws.Range("H8").Value = 0.4711538461538461
$ws.Range("I8").Value = 0.7440384615384615
$ws.Range("J8").Value = 1323
$ws.Range("K8").Value = "Content"
$ws.Range("L8").Value = "2021-09-27-03.22"
$ws.Range("M8").Value = "2021-09-27-03.22"

# row 9: Content @ 2021-10-06-19.24
$ws.Range("A9").Value = 0.7230769230769231
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0.9692307692307692
$ws.Range("D9").Value = 0.9692307692307692
$ws.Range("E9").Value = 0.9692307692307692
$ws.Range("F9").Value = 0.9692307692307692
$ws.Range("G9").Value = 0.2792307692307692
$ws.Range("H9").Value = 0.4846153846153846
$ws.Range("I9").Value = 0.7638461538461538
$ws.Range("J9").Value = 1879
$ws.Range("K9").Value = "Content"
$ws.Range("L9").Value = "2021-10-06-19.24"
$ws.Range("M9").Value = "2021-10-06-19.24"

# row 10: Files @ 2021-10-15-19.42
$ws.Range("A10").Value = 0.5
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.2475
$ws.Range("H10").Value = 0.5
$ws.Range("I10").Value = 0.7475000000000001
$ws.Range("J10").Value = 125
$ws.Range("K10").Value = "Files"
$ws.Range("L10").Value = "2021-10-15-19.42"
$ws.Range("M10").Value = "2021-10-15-19.42"

# row 11: Files @ 2021-10-15-20.05
$ws.Range("A11").Value = 0.5
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.2475
$ws.Range("H11").Value = 0.5
$ws.Range("I11").Value = 0.7475000000000001
$ws.Range("J11").Value = 125
$ws.Range("K11").Value = "Files"
$ws.Range("L11").Value = "2021-10-15-20.05"
$ws.Range("M11").Value = "2021-10-15-20.05"

# row 12: Admin @ 2021-10-18-01.33
$ws.Range("A12").Value = 0.5185185185185185
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.2505555555555555
$ws.Range("H12").Value = 0.5
$ws.Range("I12").Value = 0.7505555555555555
$ws.Range("J12").Value = 322
$ws.Range("K12").Value = "Admin"
$ws.Range("L12").Value = "2021-10-18-01.33"
$ws.Range("M12").Value = "2021-10-18-01.33"

# row 13: Mobile_App @ 2021-10-19-18.47
$ws.Range("A13").Value = 0.5411764705882353
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.2542941176470588
$ws.Range("H13").Value = 0.5
$ws.Range("I13").Value = 0.7542941176470588
$ws.Range("J13").Value = 2342
$ws.Range("K13").Value = "Mobile_App"
$ws.Range("L13").Value = "2021-10-19-18.47"
$ws.Range("M13").Value = "2021-10-19-18.47"

# row 14: Content @ 2021-10-19-18.48
$ws.Range("A14").Value = 0.6979166666666666
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0.9791666666666666
$ws.Range("D14").Value = 0.9791666666666666
$ws.Range("E14").Value = 0.9791666666666666
$ws.Range("F14").Value = 0.9791666666666666
$ws.Range("G14").Value = 0.27671875
$ws.Range("H14").Value = 0.4895833333333333
$ws.Range("I14").Value = 0.7663020833333333
$ws.Range("J14").Value = 2477
$ws.Range("K14").Value = "Content"
$ws.Range("L14").Value = "2021-10-19-18.48"
$ws.Range("M14").Value = "2021-10-19-18.48"

# row 15: User @ 2021-10-24-19.18
$ws.Range("A15").Value = 0.5
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0.9285714285714286
$ws.Range("D15").Value = 0.9285714285714286
$ws.Range("E15").Value = 0.9285714285714286
$ws.Range("F15").Value = 0.9285714285714286
$ws.Range("G15").Value = 0.2357142857142857
$ws.Range("H15").Value = 0.4642857142857143
$ws.Range("I15").Value = 0.7
$ws.Range("J15").Value = 685
$ws.Range("K15").Value = "User"
$ws.Range("L15").Value = "2021-10-24-19.18"
$ws.Range("M15").Value = "2021-10-24-19.18"

# row 16: Mobile_App @ 2021-11-03-22:59
$ws.Range("A16").Value = 0.5411764705882353
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2542941176470588
$ws.Range("H16").Value = 0.5
$ws.Range("I16").Value = 0.7542941176470588
$ws.Range("J16").Value = 2222
$ws.Range("K16").Value = "Mobile_App"
$ws.Range("L16").Value = "2021-11-03-22:59"
$ws.Range("M16").Value = "2021-11-03-22:59"

# row 17: Files @ 2021-11-03-22:59
$ws.Range("A17").Value = 0.5
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2475
$ws.Range("H17").Value = 0.5
$ws.Range("I17").Value = 0.7475000000000001
$ws.Range("J17").Value = 132
$ws.Range("K17").Value = "Files"
$ws.Range("L17").Value = "2021-11-03-22:59"
$ws.Range("M17").Value = "2021-11-03-22:59"

# row 18: Content @ 2021-11-03-22:59
$ws.Range("A18").Value = 0.71
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0.98
$ws.Range("D18").Value = 0.98
$ws.Range("E18").Value = 0.98
$ws.Range("F18").Value = 0.98
$ws.Range("G18").Value = 0.27885
$ws.Range("H18").Value = 0.49
$ws.Range("I18").Value = 0.76885
$ws.Range("J18").Value = 2647
$ws.Range("K18").Value = "Content"
$ws.Range("L18").Value = "2021-11-03-22:59"
$ws.Range("M18").Value = "2021-11-03-22:59"

# row 19: User @ 2021-11-03-22:59
$ws.Range("A19").Value = 0.5
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0.9285714285714286
$ws.Range("D19").Value = 0.9285714285714286
$ws.Range("E19").Value = 0.9285714285714286
$ws.Range("F19").Value = 0.9285714285714286
$ws.Range("G19").Value = 0.2357142857142857
$ws.Range("H19").Value = 0.4642857142857143
$ws.Range("I19").Value = 0.7
$ws.Range("J19").Value = 688
$ws.Range("K19").Value = "User"
$ws.Range("L19").Value = "2021-11-03-22:59"
$ws.Range("M19").Value = "2021-11-03-22:59"

# row 20: Admin @ 2021-11-04-00.35
$ws.Range("A20").Value = 0.6744186046511628
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.2762790697674419
$ws.Range("H20").Value = 0.5
$ws.Range("I20").Value = 0.7762790697674419
$ws.Range("J20").Value = 708
$ws.Range("K20").Value = "Admin"
$ws.Range("L20").Value = "2021-11-04-00.35"
$ws.Range("M20").Value = "2021-11-04-00.35"

# row 21: Admin @ 2021-11-08-03.34
$ws.Range("A21").Value = 0.7037037037037037
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.2811111111111111
$ws.Range("H21").Value = 0.5
$ws.Range("I21").Value = 0.7811111111111111
$ws.Range("J21").Value = 1051
$ws.Range("K21").Value = "Admin"
$ws.Range("L21").Value = "2021-11-08-03.34"
$ws.Range("M21").Value = "2021-11-08-03.34"

# row 22: Mobile_App @ 2021-11-08-20.15
$ws.Range("A22").Value = 0.5402298850574713
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 0.2541379310344828
$ws.Range("H22").Value = 0.5
$ws.Range("I22").Value = 0.7541379310344828
$ws.Range("J22").Value = 2340
$ws.Range("K22").Value = "Mobile_App"
$ws.Range("L22").Value = "2021-11-08-20.15"
$ws.Range("M22").Value = "2021-11-08-20.15"

# row 23: Content @ 2021-11-10-17.00
$ws.Range("A23").Value = 0.7032967032967034
$ws.Range("B23").Value = 0.02197802197802198
$ws.Range("C23").Value = 0.978021978021978
$ws.Range("D23").Value = 0.978021978021978
$ws.Range("E23").Value = 0.978021978021978
$ws.Range("F23").Value = 0.978021978021978
$ws.Range("G23").Value = 0.281043956043956
$ws.Range("H23").Value = 0.4890109890109889
$ws.Range("I23").Value = 0.770054945054945
$ws.Range("J23").Value = 2038
$ws.Range("K23").Value = "Content"
$ws.Range("L23").Value = "2021-11-10-17.00"
$ws.Range("M23").Value = "2021-11-10-17.00"

# row 24: Mobile_App @ 2021-11-10-17.28
$ws.Range("A24").Value = 0.5444444444444444
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 0.2548333333333334
$ws.Range("H24").Value = 0.5
$ws.Range("I24").Value = 0.7548333333333334
$ws.Range("J24").Value = 2439
$ws.Range("K24").Value = "Mobile_App"
$ws.Range("L24").Value = "2021-11-10-17.28"
$ws.Range("M24").Value = "2021-11-10-17.28"

